# Edit: re-apply the built-in "Themed Style 1 - Accent 1" table style to the
# three data tables (previously a custom-defined "Table_0" style), and switch
# the deck's design/theme colour scheme from the "Integral" (Red Violet)
# palette back to the standard Office theme palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables (slides 14, 15, 16) -----------------------
$newStyleId = "{AE9A80AF-9B54-497F-9FD2-4BB8B3F22728}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Swap the presentation theme colours back to the Office defaults -----
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in that order) expressed as
# COM BGR long values.
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
